$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8342026472091675
$ws.Range("B1").Value = 2.607484579086304
$ws.Range("C1").Value = 1.128958702087402
$ws.Range("D1").Value = 1.119112730026245
$ws.Range("E1").Value = 1.292571783065796
